$d = $word.ActiveDocument

# ============================================================
# Pembimbing 1 paragraph (first occurrence only, paragraph 18)
# Before: ...M.T.<tab><tab><tab><tab>(paraf:<tab>            )   [12 spaces, underlined]
# After : ...M.T.<tab><tab><tab>(paraf:<tab><tab><tab>)          [no underline]
# ============================================================
$p1 = $d.Paragraphs(18)

# Step 1: drop the extra tab that used to sit right before "(paraf:" --
# leave the "(paraf:" run itself untouched.
$r1a = $p1.Range.Duplicate
$ok1a = $r1a.Find.Execute( `
    "^t(paraf:", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "(paraf:", 2)

# Step 2: turn "<tab>            )" (tab + 12 spaces + close-paren) into three tabs + close-paren.
$r1b = $p1.Range.Duplicate
$ok1b = $r1b.Find.Execute( `
    "^t            )", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "^t^t^t)", 2)

# Step 3: strip the leftover underline from the newly written tail.
$r1c = $p1.Range.Duplicate
$ok1c = $r1c.Find.Execute( `
    "^t^t^t)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "", 0)
$r1c.Font.Underline = 0

# ============================================================
# Pembimbing 2 paragraph (first occurrence only, paragraph 19)
# Before: ...M.Eng.<tab><tab><tab><tab><tab>(paraf:<tab>______)   [underlined]
# After : ...M.Eng.<tab><tab><tab><tab><tab>(paraf:<tab><tab>)    [no underline]
# Keep the "(paraf:" run itself untouched -- only replace the part after it.
# ============================================================
$p2 = $d.Paragraphs(19)

# Step 1: turn "<tab>______)" into two tabs + close-paren.
$r2a = $p2.Range.Duplicate
$ok2a = $r2a.Find.Execute( `
    "^t______)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "^t^t)", 2)

# Step 2: strip the leftover underline from the newly written tail.
$r2b = $p2.Range.Duplicate
$ok2b = $r2b.Find.Execute( `
    "^t^t)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "", 0)
$r2b.Font.Underline = 0

Write-Output "Pembimbing1: step1=$ok1a step2=$ok1b step3=$ok1c"
Write-Output "Pembimbing2: step1=$ok2a step2=$ok2b"
Write-Output "P18: [$($d.Paragraphs(18).Range.Text)]"
Write-Output "P19: [$($d.Paragraphs(19).Range.Text)]"
